$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '56.730.49'
$ws.Range("E2").Value = '  +1.81%  '

# Row 3
$ws.Range("D3").Value = '2.331.30'
$ws.Range("E3").Value = '  +1.38%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$cellStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.34'
$ws.Range("D5").Style = $cellStyle
$ws.Range("E5").Value = '  +0.27%  '

# Row 6
$cellStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.57'
$ws.Range("D6").Style = $cellStyle
$ws.Range("E6").Value = '  +3.01%  '

# Row 7
$ws.Range("E7").Value = '  +0.38%  '

# Row 8
$cellStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.537'
$ws.Range("D8").Style = $cellStyle
$ws.Range("E8").Value = '  +0.78%  '

# Row 9
$ws.Range("D9").Value = '2.336.65'
$ws.Range("E9").Value = '  +0.76%  '

# Row 10
$ws.Range("E10").Value = '  -0.50%  '

# Row 11
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cellStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.37'
$ws.Range("D11").Style = $cellStyle
$ws.Range("E11").Value = '  +5.32%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cellStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.152'
$ws.Range("D12").Style = $cellStyle
$ws.Range("E12").Value = '  -1.67%  '

# Row 13
$ws.Range("E13").Value = '  -0.17%  '

# Row 14
$cellStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.91'
$ws.Range("D14").Style = $cellStyle
$ws.Range("E14").Value = '  -0.32%  '

# Row 15
$ws.Range("D15").Value = '2.744.38'
$ws.Range("E15").Value = '  +1.16%  '

# Row 16
$ws.Range("D16").Value = '56.686.82'
$ws.Range("E16").Value = '  +1.39%  '

# Row 17
$ws.Range("E17").Value = '  +0.14%  '

# Row 18
$ws.Range("D18").Value = '2.339.47'
$ws.Range("E18").Value = '  +2.15%  '

# Row 19
$cellStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.53'
$ws.Range("D19").Style = $cellStyle
$ws.Range("E19").Value = '  -0.10%  '

# Row 20
$cellStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '326.25'
$ws.Range("D20").Style = $cellStyle
$ws.Range("E20").Value = '  +2.22%  '

# Row 21
$ws.Range("E21").Value = '  -0.36%  '

# Row 22
$ws.Range("E22").Value = '  +0.37%  '

# Row 23
$ws.Range("E23").Value = '  +0.41%  '

# Row 24
$cellStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.89'
$ws.Range("D24").Style = $cellStyle
$ws.Range("E24").Value = '  +0.87%  '

# Row 25
$ws.Range("E25").Value = '  +4.95%  '

# Row 26
$ws.Range("E26").Value = '  +0.69%  '

# Row 27
$ws.Range("E27").Value = '  +4.71%  '

# Row 28
$cellStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.30'
$ws.Range("D28").Style = $cellStyle
$ws.Range("E28").Value = '  +9.91%  '

# Row 29
$cellStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.05'
$ws.Range("D29").Style = $cellStyle
$ws.Range("E29").Value = '  -0.40%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0739'
$ws.Range("E30").Value = '  +2.64%  '

# Row 31
$ws.Range("E31").Value = '  +1.41%  '

# Row 32
$cellStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.21'
$ws.Range("D32").Style = $cellStyle
$ws.Range("E32").Value = '  -0.29%  '

# Row 33
$cellStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.50'
$ws.Range("D33").Style = $cellStyle
$ws.Range("E33").Value = '  +1.67%  '

# Row 34
$ws.Range("E34").Value = '  +0.02%  '

# Row 35
$cellStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = $cellStyle
$ws.Range("E35").Value = '  +0.67%  '

# Row 36
$cellStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.27'
$ws.Range("D36").Style = $cellStyle
$ws.Range("E36").Value = '  +1.10%  '

# Row 37
$ws.Range("E37").Value = '  -0.87%  '

# Row 38
$cellStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.00'
$ws.Range("D38").Style = $cellStyle
$ws.Range("E38").Value = '  +1.07%  '

# Row 39
$ws.Range("E39").Value = '  +3.35%  '

# Row 40
$cellStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.30'
$ws.Range("D40").Style = $cellStyle
$ws.Range("E40").Value = '  +2.95%  '

# Row 41
$cellStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '142.41'
$ws.Range("D41").Style = $cellStyle
$ws.Range("E41").Value = '  +3.24%  '

# Row 42
$ws.Range("E42").Value = '  -0.12%  '

# Row 43
$cellStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.61'
$ws.Range("D43").Style = $cellStyle
$ws.Range("E43").Value = '  +1.49%  '

# Row 44
$cellStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '277.87'
$ws.Range("D44").Style = $cellStyle
$ws.Range("E44").Value = '  +5.29%  '

# Row 45
$cellStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.15'
$ws.Range("D45").Style = $cellStyle
$ws.Range("E45").Value = '  +0.63%  '

# Row 46
$ws.Range("E46").Value = '  +1.22%  '

# Row 47
$cellStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0506'
$ws.Range("D47").Style = $cellStyle
$ws.Range("E47").Value = '  -0.79%  '

# Row 48
$cellStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.562'
$ws.Range("D48").Style = $cellStyle
$ws.Range("E48").Value = '  +1.38%  '

# Row 49
$cellStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0219'
$ws.Range("D49").Style = $cellStyle
$ws.Range("E49").Value = '  +1.61%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cellStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.02'
$ws.Range("D50").Style = $cellStyle
$ws.Range("E50").Value = '  +7.87%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cellStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.58'
$ws.Range("D51").Style = $cellStyle
$ws.Range("E51").Value = '  +4.84%  '
